$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1, copying the style from H1 (bold + border + centered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill I and J columns with data for rows 2-38
$data = @(
    ,@(2, 8, 8)
    ,@(3, 6, 6)
    ,@(4, 7, 7)
    ,@(5, 7, 7)
    ,@(6, 7, 7)
    ,@(7, 6, 6)
    ,@(8, 6, 6)
    ,@(9, 7, 7)
    ,@(10, 7, 7)
    ,@(11, 6, 7)
    ,@(12, 7, 7)
    ,@(13, 6, 6)
    ,@(14, 5, 6)
    ,@(15, 7, 7)
    ,@(16, 8, 8)
    ,@(17, 7, 7)
    ,@(18, 6, 6)
    ,@(19, 6, 6)
    ,@(20, 8, 8)
    ,@(21, 6, 7)
    ,@(22, 6, 6)
    ,@(23, 5, 5)
    ,@(24, 5, 5)
    ,@(25, 8, 8)
    ,@(26, 7, 7)
    ,@(27, 6, 7)
    ,@(28, 8, 8)
    ,@(29, 8, 8)
    ,@(30, 6, 6)
    ,@(31, 8, 8)
    ,@(32, 9, 9)
    ,@(33, 6, 6)
    ,@(34, 9, 9)
    ,@(35, 8, 8)
    ,@(36, 8, 8)
    ,@(37, 5, 5)
    ,@(38, 5, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
